{"js": "// Update the division exercises in the table to the new problem set.\n// The mapping below is derived directly from the source OOXML diff\n// (old expression text -> new expression text); every value in the table\n// appears exactly once, so an exact, case-sensitive search is unambiguous.\nconst replacements = [\n  [\"152\u00f77=\", \"228\u00f79=\"],\n  [\"505\u00f79=\", \"158\u00f77=\"],\n  [\"992\u00f77=\", \"831\u00f79=\"],\n  [\"123\u00f76=\", \"334\u00f77=\"],\n  [\"219\u00f74=\", \"739\u00f74=\"],\n  [\"516\u00f78=\", \"970\u00f76=\"],\n  [\"234\u00f74=\", \"182\u00f76=\"],\n  [\"238\u00f79=\", \"535\u00f75=\"],\n  [\"667\u00f74=\", \"315\u00f74=\"],\n  [\"757\u00f77=\", \"870\u00f78=\"],\n  [\"283\u00f73=\", \"150\u00f73=\"],\n  [\"103\u00f79=\", \"223\u00f75=\"],\n  [\"733\u00f72=\", \"619\u00f75=\"],\n  [\"778\u00f73=\", \"755\u00f74=\"],\n  [\"865\u00f75=\", \"423\u00f79=\"],\n  [\"943\u00f79=\", \"902\u00f79=\"],\n  [\"587\u00f73=\", \"496\u00f75=\"],\n  [\"290\u00f72=\", \"671\u00f79=\"],\n  [\"461\u00f79=\", \"290\u00f78=\"],\n  [\"692\u00f73=\", \"673\u00f78=\"],\n  [\"106\u00f75=\", \"271\u00f79=\"],\n  [\"397\u00f76=\", \"564\u00f72=\"],\n  [\"839\u00f73=\", \"918\u00f75=\"],\n  [\"643\u00f77=\", \"324\u00f79=\"],\n  [\"778\u00f75=\", \"895\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Expected to find \"${oldText}\" in the document, but it was not found.`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division exercises in the table to the new problem set.\n# The mapping below is derived directly from the source OOXML diff\n# (old expression text -> new expression text); every value in the table\n# appears exactly once, so an exact Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"152\u00f77=\", \"228\u00f79=\"),\n    @(\"505\u00f79=\", \"158\u00f77=\"),\n    @(\"992\u00f77=\", \"831\u00f79=\"),\n    @(\"123\u00f76=\", \"334\u00f77=\"),\n    @(\"219\u00f74=\", \"739\u00f74=\"),\n    @(\"516\u00f78=\", \"970\u00f76=\"),\n    @(\"234\u00f74=\", \"182\u00f76=\"),\n    @(\"238\u00f79=\", \"535\u00f75=\"),\n    @(\"667\u00f74=\", \"315\u00f74=\"),\n    @(\"757\u00f77=\", \"870\u00f78=\"),\n    @(\"283\u00f73=\", \"150\u00f73=\"),\n    @(\"103\u00f79=\", \"223\u00f75=\"),\n    @(\"733\u00f72=\", \"619\u00f75=\"),\n    @(\"778\u00f73=\", \"755\u00f74=\"),\n    @(\"865\u00f75=\", \"423\u00f79=\"),\n    @(\"943\u00f79=\", \"902\u00f79=\"),\n    @(\"587\u00f73=\", \"496\u00f75=\"),\n    @(\"290\u00f72=\", \"671\u00f79=\"),\n    @(\"461\u00f79=\", \"290\u00f78=\"),\n    @(\"692\u00f73=\", \"673\u00f78=\"),\n    @(\"106\u00f75=\", \"271\u00f79=\"),\n    @(\"397\u00f76=\", \"564\u00f72=\"),\n    @(\"839\u00f73=\", \"918\u00f75=\"),\n    @(\"643\u00f77=\", \"324\u00f79=\"),\n    @(\"778\u00f75=\", \"895\u00f73=\"),\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $find.Replacement.Text,\n        $wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Expected to find $oldText in the document, but it was not found.\"\n    }\n}\n"}
